# Apply updates described by the commit:
#  - Added transport, commercial and industrial demands
#  - Added factor for energy losses
#  - Merged all regions for demand and supply

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Battery (row 2) charging/discharging efficiency factors
$ws.Range("C2").Value = 0.9
$ws.Range("D2").Value = 0.9

# Update Hydro (row 4) charging efficiency and C-rate (energy loss factor)
$ws.Range("C4").Value = 0.9
$ws.Range("E4").Value = 0.001

# Update the active selection to match the author's saved view state
$ws.Range("E3").Select()
